# db/PATIENT.xlsx changes:
#  - E3 was stored as text "1234546798"; convert it to a real number.
#  - Append a new patient row (row 4) with the given fields.
#    E4 ("0987654322") has a leading zero so it must stay text, not be
#    coerced into a number that would drop the leading digit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: E3 becomes a numeric value instead of inline text
$ws.Range("E3").Value = 1234546798

# Row 4: new patient record
$ws.Range("A4").Value = "abhishek@hello.com"
$ws.Range("B4").Value = "password"
$ws.Range("C4").Value = "MALE"
$ws.Range("D4").Value = 32

# Force E4 to remain text so the leading zero in the phone number survives,
# then restore the default style so no stray number-format style sticks.
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0987654322"
$ws.Range("E4").Style = "Normal"
